$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G3").Value = 3000
$ws.Range("G4").Value = 1000
$ws.Range("G10").Value = 8000
$ws.Range("G12").Value = 3000
$ws.Range("G14").Value = 0
$ws.Range("G16").Value = 3000
$ws.Range("G19").Value = 1000
$ws.Range("G39").Value = 0
$ws.Range("G40").Value = 2000
$ws.Range("G44").Value = 2000
$ws.Range("G45").Value = 4500
$ws.Range("G49").Value = 2000
$ws.Range("G50").Value = 0
$ws.Range("G51").Value = 3500
$ws.Range("G53").Value = 0
$ws.Range("G55").Value = 1000
$ws.Range("G56").Value = 5000
$ws.Range("G58").Value = 57000
